$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value2 = "360+5"
$ws.Range("E4").Value2 = "264+5"
$ws.Range("E5").Value2 = "288+5"
$ws.Range("E6").Value2 = "312+5"
$ws.Range("E7").Value2 = "360+5"
$ws.Range("E8").Value2 = "432+5"
$ws.Range("E9").Value2 = "312+5"
$ws.Range("E10").Value2 = "336+5"
